$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.930.71'
$ws.Range("E2").Value = '  -3.08%  '

$ws.Range("D3").Value = '3.838.03'
$ws.Range("E3").Value = '  -2.46%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.20'
$ws.Range("E5").Value = '  -1.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.97'
$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("D7").Value = '3.837.06'
$ws.Range("E7").Value = '  -2.45%  '

$ws.Range("E9").Value = '  -1.94%  '

$ws.Range("E10").Value = '  -3.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.47'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -2.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000260'
$ws.Range("E13").Value = '  +1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.05'
$ws.Range("E14").Value = '  -3.89%  '

$ws.Range("D15").Value = '4.476.77'
$ws.Range("E15").Value = '  -2.57%  '

$ws.Range("D16").Value = '3.836.66'
$ws.Range("E16").Value = '  -1.32%  '

$ws.Range("D17").Value = '67.925.61'
$ws.Range("E17").Value = '  -3.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.31'
$ws.Range("E18").Value = '  -2.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.41'
$ws.Range("E19").Value = '  -3.25%  '

$ws.Range("E20").Value = '  -1.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.07'
$ws.Range("E21").Value = '  -0.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.82'
$ws.Range("E22").Value = '  -5.79%  '

$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("E24").Value = '  -3.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.88'
$ws.Range("E25").Value = '  -3.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("E26").Value = '  -3.10%  '

$ws.Range("E27").Value = '  -2.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("E28").Value = '  -1.54%  '

$ws.Range("E30").Value = '  -1.55%  '

$ws.Range("D31").Value = '3.983.18'
$ws.Range("E31").Value = '  -2.52%  '

$ws.Range("E32").Value = '  -2.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.38'
$ws.Range("E33").Value = '  -3.02%  '

$ws.Range("E34").Value = '  -6.15%  '

$ws.Range("E35").Value = '  -0.70%  '

$ws.Range("D36").Value = '3.796.96'
$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("E37").Value = '  -3.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.65'
$ws.Range("E38").Value = '  +9.81%  '

$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("E39").Value = '  -2.45%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.140'
$ws.Range("E40").Value = '  -1.11%  '

$ws.Range("E41").Value = '  -4.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.315'
$ws.Range("E43").Value = '  -4.81%  '

$ws.Range("E44").Value = '  -7.27%  '

$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.73'
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000297'
$ws.Range("E46").Value = '  +6.23%  '

$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '419.41'
$ws.Range("E47").Value = '  -4.75%  '

$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.72'
$ws.Range("E49").Value = '  -3.46%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.58'
$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.27'
$ws.Range("E51").Value = '  +3.99%  '
